# Format updates for cash flow currency
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "% Complete" values in column E ---
$ws.Range("E12").Value = 0.9
$ws.Range("E13").Value = 0.9
$ws.Range("E14").Value = 0.9
$ws.Range("E16").Value = 0.9

# --- Re-format a few progress cells to match the "Actual (beyond plan)" fill ---
# Cells I12:M12 already carry the target format; copy it onto N12/O12.
$ws.Range("I12").Copy()
$ws.Range("N12:O12").PasteSpecial(-4122)

$ws.Range("I13").Copy()
$ws.Range("N13:O13").PasteSpecial(-4122)

$ws.Range("I16").Copy()
$ws.Range("N14").PasteSpecial(-4122)

[void]($excel.CutCopyMode = 0)

# --- Update the active selection shown when the sheet is reopened ---
[void]$ws.Range("E17").Select()
